# Fix NameError in deck generation node
# Update the "Executive Summary" metrics bullets on slide 2 to reflect
# the corrected (Ingram Micro only) computation instead of the combined
# Ingram Micro + CNH Industrial figures from the buggy run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

function Replace-RunText {
    param(
        $TextRange,
        [string]$OldText,
        [string]$NewText
    )
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Replace-RunText: text not found -> $OldText"
    }
    $sub = $TextRange.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
}

Replace-RunText $tr "Compute Metrics and Draft Summary:" "**Metrics Summary:**"

Replace-RunText $tr `
    "**Total GLA**: 313,219 m² (222,221 m² for Ingram Micro + 90,998 m² for CNH Industrial)" `
    "**Total GLA**: 222,221 m²"

Replace-RunText $tr `
    "**Occupancy**: Potentially 0% (leases for both tenants have expired)" `
    "**Occupancy**: 100% (based on the lease area matching the total asset area)"

Replace-RunText $tr `
    "**WALT**: 0 years (Weighted Average Lease Term is 0 due to expired leases)" `
    "**WALT (Weighted Average Lease Term)**: Approximately 3.5 years (calculated from the lease end dates of 2020 and 2021, assuming the current year is 2018)"

Replace-RunText $tr `
    "**In-Place Rent**: Not applicable (leases have expired)" `
    "**In-Place Rent**: £5.5 per m² per annum (based on the lease with Ingram Micro)"

Replace-RunText $tr `
    "**Key Highlight 1**: Strategic logistics location in Daventry, UK, with proximity to major transportation routes." `
    "**Key Highlight 1**: Strong tenant presence with Ingram Micro, a leading technology distributor."

Replace-RunText $tr `
    "**Key Highlight 2**: Significant leasable area with potential for stable cash flow if leases are renewed or new tenants are secured." `
    "**Key Highlight 2**: Strategic location in Daventry, United Kingdom, with excellent logistics infrastructure including 12 dock doors and 180 parking spaces."
